$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing row 15 ("REMESAS") is being split into three separate rows,
# one per underlying component variable (REMESAS1/REMESAS2/REMESAS3).
# Insert two new rows right after row 15 to make room for REMESAS2/REMESAS3,
# which pushes the old rows 16-18 down to 18-20.
$ws.Rows("16:17").Insert()

# Row 15: REMESAS -> REMESAS1, and the H:Q columns move from the combined
# average formula text to the first individual component code.
$ws.Range("A15").Value = "REMESAS1"
$ws.Range("B15").Value = "REMESAS1"
$ws.Range("H15:Q15").Value = "P05A17B"

# Row 16 (new): REMESAS2 - only the H:Q year columns are populated.
$ws.Range("A16").Value = "REMESAS2"
$ws.Range("B16").Value = "REMESAS2"
$ws.Range("H16:Q16").Value = "P05A18B"

# Row 17 (new): REMESAS3 - only the H:Q year columns are populated.
$ws.Range("A17").Value = "REMESAS3"
$ws.Range("B17").Value = "REMESAS3"
$ws.Range("H17:Q17").Value = "P05A19B"

# Match the final cursor position recorded in the saved workbook.
[void]$ws.Range("J12").Select()
